$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

function Add-EmptyCell {
    param($cell)
    $cell.NumberFormat = "@"
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("A2") "03/09/2021"
Set-TextValue $ws.Range("G2") "0.660000"
Set-TextValue $ws.Range("H2") "29101.23"
Set-TextValue $ws.Range("J2") "121588240.00"
Add-EmptyCell $ws.Range("M2")

# Row 3
Set-TextValue $ws.Range("A3") "03/09/2021"
Set-TextValue $ws.Range("G3") "20.460000"
Set-TextValue $ws.Range("H3") "10004917.40"
Set-TextValue $ws.Range("I3") "8.23%"
Set-TextValue $ws.Range("J3") "121588240.00"
Add-EmptyCell $ws.Range("M3")

# Row 4
Set-TextValue $ws.Range("A4") "03/09/2021"
Set-TextValue $ws.Range("G4") "71.460000"
Set-TextValue $ws.Range("H4") "1106843.94"
Set-TextValue $ws.Range("I4") "0.91%"
Set-TextValue $ws.Range("J4") "121588240.00"
Add-EmptyCell $ws.Range("M4")

# Row 5
Set-TextValue $ws.Range("A5") "03/09/2021"
Set-TextValue $ws.Range("G5") "18.110000"
Set-TextValue $ws.Range("H5") "3742757.48"
Set-TextValue $ws.Range("I5") "3.08%"
Set-TextValue $ws.Range("J5") "121588240.00"
Add-EmptyCell $ws.Range("M5")

# Row 6
Set-TextValue $ws.Range("A6") "03/09/2021"
Set-TextValue $ws.Range("G6") "30.860000"
Set-TextValue $ws.Range("H6") "11647304.64"
Set-TextValue $ws.Range("J6") "121588240.00"
Add-EmptyCell $ws.Range("M6")

# Row 7
Set-TextValue $ws.Range("A7") "03/09/2021"
Set-TextValue $ws.Range("G7") "11.660000"
Set-TextValue $ws.Range("H7") "3450951.90"
Set-TextValue $ws.Range("I7") "2.84%"
Set-TextValue $ws.Range("J7") "121588240.00"
Add-EmptyCell $ws.Range("M7")

# Row 8
Set-TextValue $ws.Range("A8") "03/09/2021"
Set-TextValue $ws.Range("G8") "9.440000"
Set-TextValue $ws.Range("H8") "5588281.76"
Set-TextValue $ws.Range("I8") "4.60%"
Set-TextValue $ws.Range("J8") "121588240.00"
Add-EmptyCell $ws.Range("M8")

# Row 9
Set-TextValue $ws.Range("A9") "03/09/2021"
Set-TextValue $ws.Range("G9") "5.510000"
Set-TextValue $ws.Range("H9") "3528825.51"
Set-TextValue $ws.Range("I9") "2.90%"
Set-TextValue $ws.Range("J9") "121588240.00"
Add-EmptyCell $ws.Range("M9")

# Row 10
Set-TextValue $ws.Range("A10") "03/09/2021"
Set-TextValue $ws.Range("G10") "5.010000"
Set-TextValue $ws.Range("H10") "2132882.25"
Set-TextValue $ws.Range("I10") "1.75%"
Set-TextValue $ws.Range("J10") "121588240.00"
Add-EmptyCell $ws.Range("M10")

# Row 11
Set-TextValue $ws.Range("A11") "03/09/2021"
Set-TextValue $ws.Range("G11") "37.790000"
Set-TextValue $ws.Range("H11") "3664987.57"
Set-TextValue $ws.Range("I11") "3.01%"
Set-TextValue $ws.Range("J11") "121588240.00"
Add-EmptyCell $ws.Range("M11")

# Row 12
Set-TextValue $ws.Range("A12") "03/09/2021"
Set-TextValue $ws.Range("G12") "214.720000"
Set-TextValue $ws.Range("H12") "14455809.28"
Set-TextValue $ws.Range("I12") "11.89%"
Set-TextValue $ws.Range("J12") "121588240.00"
Add-EmptyCell $ws.Range("M12")

# Row 13
Set-TextValue $ws.Range("A13") "03/09/2021"
Set-TextValue $ws.Range("G13") "47.420000"
Set-TextValue $ws.Range("H13") "5579200.10"
Set-TextValue $ws.Range("I13") "4.59%"
Set-TextValue $ws.Range("J13") "121588240.00"
Add-EmptyCell $ws.Range("M13")

# Row 14
Set-TextValue $ws.Range("A14") "03/09/2021"
Set-TextValue $ws.Range("G14") "165.650000"
Set-TextValue $ws.Range("H14") "4143734.75"
Set-TextValue $ws.Range("I14") "3.41%"
Set-TextValue $ws.Range("J14") "121588240.00"
Add-EmptyCell $ws.Range("M14")

# Row 15
Set-TextValue $ws.Range("A15") "03/09/2021"
Set-TextValue $ws.Range("C15") "IMCC CN"
Set-TextValue $ws.Range("G15") "10.200000"
Set-TextValue $ws.Range("H15") "40233.51"
Set-TextValue $ws.Range("J15") "121588240.00"
Add-EmptyCell $ws.Range("M15")

# Row 16
Set-TextValue $ws.Range("A16") "03/09/2021"
Set-TextValue $ws.Range("G16") "4.780000"
Set-TextValue $ws.Range("H16") "878640.48"
Set-TextValue $ws.Range("I16") "0.72%"
Set-TextValue $ws.Range("J16") "121588240.00"
Add-EmptyCell $ws.Range("M16")

# Row 17
Set-TextValue $ws.Range("A17") "03/09/2021"
Set-TextValue $ws.Range("G17") "0.475000"
Set-TextValue $ws.Range("H17") "1273962.67"
Set-TextValue $ws.Range("I17") "1.05%"
Set-TextValue $ws.Range("J17") "121588240.00"
Add-EmptyCell $ws.Range("M17")

# Row 18
Set-TextValue $ws.Range("A18") "03/09/2021"
Set-TextValue $ws.Range("G18") "1.390000"
Set-TextValue $ws.Range("H18") "1687404.40"
Set-TextValue $ws.Range("J18") "121588240.00"
Add-EmptyCell $ws.Range("M18")

# Row 19
Set-TextValue $ws.Range("A19") "03/09/2021"
Set-TextValue $ws.Range("G19") "2.830000"
Set-TextValue $ws.Range("H19") "2468218.46"
Set-TextValue $ws.Range("I19") "2.03%"
Set-TextValue $ws.Range("J19") "121588240.00"
Add-EmptyCell $ws.Range("M19")

# Row 20
Set-TextValue $ws.Range("A20") "03/09/2021"
Set-TextValue $ws.Range("G20") "121.560000"
Set-TextValue $ws.Range("H20") "4337382.36"
Set-TextValue $ws.Range("I20") "3.57%"
Set-TextValue $ws.Range("J20") "121588240.00"
Add-EmptyCell $ws.Range("M20")

# Row 21
Set-TextValue $ws.Range("A21") "03/09/2021"
Set-TextValue $ws.Range("G21") "2.590000"
Set-TextValue $ws.Range("H21") "6581803.79"
Set-TextValue $ws.Range("I21") "5.41%"
Set-TextValue $ws.Range("J21") "121588240.00"
Add-EmptyCell $ws.Range("M21")

# Row 22
Set-TextValue $ws.Range("A22") "03/09/2021"
Set-TextValue $ws.Range("G22") "17.640000"
Set-TextValue $ws.Range("H22") "8932896.00"
Set-TextValue $ws.Range("I22") "7.35%"
Set-TextValue $ws.Range("J22") "121588240.00"
Add-EmptyCell $ws.Range("M22")

# Row 23
Set-TextValue $ws.Range("A23") "03/09/2021"
Set-TextValue $ws.Range("G23") "22.420000"
Set-TextValue $ws.Range("H23") "1143958.08"
Set-TextValue $ws.Range("I23") "0.94%"
Set-TextValue $ws.Range("J23") "121588240.00"
Add-EmptyCell $ws.Range("M23")

# Row 24
Set-TextValue $ws.Range("A24") "03/09/2021"
Set-TextValue $ws.Range("G24") "13.750000"
Set-TextValue $ws.Range("H24") "10482890.00"
Set-TextValue $ws.Range("I24") "8.62%"
Set-TextValue $ws.Range("J24") "121588240.00"
Add-EmptyCell $ws.Range("M24")

# Row 25
Set-TextValue $ws.Range("A25") "03/09/2021"
Set-TextValue $ws.Range("C25") "VLNS CN"
Set-TextValue $ws.Range("G25") "1.700000"
Set-TextValue $ws.Range("H25") "53644.68"
Set-TextValue $ws.Range("J25") "121588240.00"
Add-EmptyCell $ws.Range("M25")

# Row 26
Set-TextValue $ws.Range("A26") "03/09/2021"
Set-TextValue $ws.Range("H26") "589368.38"
Set-TextValue $ws.Range("J26") "121588240.00"
Add-EmptyCell $ws.Range("M26")

# Row 27
Set-TextValue $ws.Range("A27") "03/09/2021"
Set-TextValue $ws.Range("G27") "3.960000"
Set-TextValue $ws.Range("H27") "1823877.00"
Set-TextValue $ws.Range("I27") "1.50%"
Set-TextValue $ws.Range("J27") "121588240.00"
Add-EmptyCell $ws.Range("M27")

# Row 28
Set-TextValue $ws.Range("A28") "03/09/2021"
Set-TextValue $ws.Range("G28") "4.140000"
Set-TextValue $ws.Range("H28") "1126013.76"
Set-TextValue $ws.Range("I28") "0.93%"
Set-TextValue $ws.Range("J28") "121588240.00"
Add-EmptyCell $ws.Range("M28")
